$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: row -> [D, M, N, O, P, S] target values (per commit: weekly fruit/vegetable refresh)
$rowsData = @{
    2 = @(44467,200,20000,21000,20500,1025)
    3 = @(44410,200,20000,21000,20500,1025)
    4 = @(44431,160,21000,22000,21500,1075)
    5 = @(44418,200,20000,21000,20500,1025)
    6 = @(44781,160,23000,24000,23500,1175)
    7 = @(44809,60,27000,28000,27500,1375)
    8 = @(44474,200,19000,20000,19500,975)
    9 = @(44810,100,27000,28000,27500,1375)
    10 = @(44466,100,20000,21000,20500,1025)
    11 = @(44407,160,20000,21000,20500,1025)
    12 = @(44350,160,19000,20000,19500,975)
    13 = @(44473,40,19500,20000,19750,988)
    14 = @(44301,100,18000,19000,18500,925)
    15 = @(44326,160,19500,20000,19750,988)
    16 = @(44879,100,28000,30000,29000,1450)
    17 = @(44434,100,20000,21000,20500,1025)
    18 = @(44784,160,27000,28000,27500,1375)
    19 = @(44417,160,20000,21000,20500,1025)
    20 = @(44315,100,20000,21000,20500,1025)
    21 = @(44428,100,20000,21000,20500,1025)
    22 = @(44445,160,20000,21000,20500,1025)
    23 = @(44343,100,19500,20000,19750,988)
    24 = @(44333,100,19500,20000,19750,988)
    25 = @(44427,200,20000,21000,20500,1025)
    26 = @(44882,120,28000,30000,29000,1450)
    27 = @(44335,200,19000,20000,19500,975)
    28 = @(44336,100,19500,20000,19750,988)
    29 = @(44782,200,23500,24000,23750,1188)
    30 = @(44776,160,23000,24000,23500,1175)
    31 = @(44435,260,20000,22000,21115,1056)
    32 = @(44442,140,20000,21000,20500,1025)
    33 = @(44364,140,20000,21000,20500,1025)
    34 = @(44448,100,20000,21000,20500,1025)
    35 = @(44420,160,20000,21000,20500,1025)
    36 = @(44441,160,20000,21000,20500,1025)
    37 = @(44880,100,28000,30000,29000,1450)
    38 = @(44778,100,23000,24000,23500,1175)
    39 = @(44462,100,19500,20000,19750,988)
    40 = @(44874,240,29000,30000,29500,1475)
    41 = @(44365,100,20000,21000,20500,1025)
}

foreach ($row in $rowsData.Keys) {
    $vals = $rowsData[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]  # D: Fecha
    $ws.Cells.Item($row, 13).Value = $vals[1]  # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals[2]  # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals[3]  # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals[4]  # P: Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals[5]  # S: Precio $/Kg
}
